# Generate Report for Archive
#
# The localization run moved from "Ready for handoff" to "In Translation"
# for this file, across the Overview sheet (one column per locale) and
# each per-locale detail sheet (Status column). Updating the status text
# also shrinks the affected status columns, since they no longer need to
# fit the longer "Ready for handoff" label.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn status is column E, de-de status is column F ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# Narrow columns E and F to fit the shorter status text. ColumnWidth only
# lands on a 1/6-character grid, so 12.5 is the input that snaps closest
# to the target (~13.41 raw units -> 13.33 on that grid).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-locale detail sheets: Status is column C ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
